# #815 - Phase accordion
#
# The "NumberOfDays" default used to seed new Events was hard-coded to 1;
# with the phase-accordion work it should default to 0 instead (days are
# now driven from the phase itself). Reset every existing template row's
# NumberOfDays (column I) on the "Events" sheet from 1 -> 0, leaving any
# row that already carries a real, non-default value (e.g. row 5, which
# holds 30 days for the Comment Period) untouched.

$wb = $excel.ActiveWorkbook
$events = $wb.Worksheets.Item("Events")

$lastRow = $events.Cells.Item($events.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $events.Cells.Item($r, 9)
    if ($cell.Value2 -eq 1) {
        $cell.Value = 0
    }
}

# Leave the workbook with the "Events" sheet frontmost/selected, with the
# cursor resting on J62 - where the review of the accordion changes ended.
$events.Activate()
$events.Range("J62").Select()
